$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 10041802
$ws.Range("C3").Value = "RAM_MPP.doc"

$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = (Get-Date -Year 2018 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0)

$ws.Range("G15").Select()
